$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Notes" header in E1 (match the style of the other header cells) ---
$ws.Range("E1").Value = "Notes"
$ws.Range("E1").Font.Name = "Calibri"
$ws.Range("E1").Font.Size = 11
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Font.Color = $ws.Range("D1").Font.Color
$ws.Range("E1").Borders.Item(9).LineStyle = 1
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4160

# --- New EV "NOS0" node rows (2030 + 2040), same shape as the other nodes ---
$ws.Range("A10").Value = "NOS0"
$ws.Range("B10").Value = "Distributed Energy"
$ws.Range("C10").Value = 2030
$ws.Range("D10").Value = 30000

$ws.Range("A11").Value = "NOS0"
$ws.Range("B11").Value = "Distributed Energy"
$ws.Range("C11").Value = 2040
$ws.Range("D11").Value = 50000

# Match the formatting used by the existing data rows
$ws.Range("B10").Font = $ws.Range("B2").Font
$ws.Range("B11").Font = $ws.Range("B2").Font

# Restore the selection that was active when the workbook was saved
$null = $ws.Range("F8").Select()
